$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.243.13'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '''1.828.82'
$ws.Range('E3').Value = '  -0.75%  '
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').Value = '''234.60'
$ws.Range('E5').Value = '  -2.02%  '
$ws.Range('D6').Value = '''0.6007'
$ws.Range('E6').Value = '  -3.98%  '
$ws.Range('D7').Value = '''1.004'
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').Value = '''0.07025'
$ws.Range('E8').Value = '  -5.50%  '
$ws.Range('D9').Value = '''0.2782'
$ws.Range('E9').Value = '  -3.92%  '
$ws.Range('D10').Value = '''23.41'
$ws.Range('E10').Value = '  -6.23%  '
$ws.Range('D11').Value = '''0.07649'
$ws.Range('D12').Value = '''1.826.92'
$ws.Range('E12').Value = '  -0.82%  '
$ws.Range('D13').Value = '''4.784'
$ws.Range('E13').Value = '  -3.90%  '
$ws.Range('D14').Value = '''0.000009934'
$ws.Range('E14').Value = '  -3.31%  '
$ws.Range('D15').Value = '''0.6256'
$ws.Range('E15').Value = '  -7.35%  '
$ws.Range('D16').Value = '''78.86'
$ws.Range('E16').Value = '  -3.72%  '
$ws.Range('D17').Value = '''29.237.25'
$ws.Range('E17').Value = '  -0.55%  '
$ws.Range('D18').Value = '''5.829'
$ws.Range('E18').Value = '  -6.44%  '
$ws.Range('D19').Value = '''223.37'
$ws.Range('E19').Value = '  -4.30%  '
$ws.Range('D20').Value = '''1.004'
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('D21').Value = '''11.67'
$ws.Range('E21').Value = '  -5.30%  '
$ws.Range('D22').Value = '''6.994'
$ws.Range('E22').Value = '  -4.44%  '
$ws.Range('D23').Value = '''1.004'
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('D24').Value = '''155.25'
$ws.Range('E24').Value = '  -2.07%  '
$ws.Range('D25').Value = '''7.959'
$ws.Range('E25').Value = '  -6.19%  '
$ws.Range('D26').Value = '''0.1293'
$ws.Range('E26').Value = '  -4.19%  '
$ws.Range('D27').Value = '''16.52'
$ws.Range('E27').Value = '  -4.67%  '
$ws.Range('D28').Value = '''1.481'
$ws.Range('E28').Value = '  +1.27%  '
$ws.Range('D29').Value = '''0.06212'
$ws.Range('E29').Value = '  -14.43%  '
$ws.Range('D30').Value = '''1.444'
$ws.Range('E30').Value = '  -2.45%  '
$ws.Range('D31').Value = '''3.834'
$ws.Range('E31').Value = '  -4.88%  '
$ws.Range('D32').Value = '''3.793'
$ws.Range('E32').Value = '  -6.61%  '
$ws.Range('D33').Value = '''1.108'
$ws.Range('E33').Value = '  -2.83%  '
$ws.Range('D34').Value = '''1.737'
$ws.Range('E34').Value = '  -4.75%  '
$ws.Range('D35').Value = '''0.6443'
$ws.Range('E35').Value = '  -7.84%  '
$ws.Range('D36').Value = '''2.547'
$ws.Range('E36').Value = '  -1.11%  '
$ws.Range('D37').Value = '''1.222.70'
$ws.Range('E37').Value = '  -0.95%  '
$ws.Range('D38').Value = '''2.739'
$ws.Range('E38').Value = '  -2.77%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '''0.01729'
$ws.Range('E39').Value = '  -6.03%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '''6.497'
$ws.Range('E40').Value = '  -6.26%  '
$ws.Range('D41').Value = '''0.8984'
$ws.Range('E41').Value = '  -5.82%  '
$ws.Range('D42').Value = '''1.004'
$ws.Range('E42').Value = '  +0.30%  '
$ws.Range('D43').Value = '''1.984.07'
$ws.Range('E43').Value = '  -1.10%  '
$ws.Range('D44').Value = '''100.30'
$ws.Range('E44').Value = '  -0.75%  '
$ws.Range('D45').Value = '''62.36'
$ws.Range('E45').Value = '  -4.82%  '
$ws.Range('D46').Value = '''0.00000000115'
$ws.Range('E46').Value = '  -1.81%  '
$ws.Range('D47').Value = '''8.543'
$ws.Range('E47').Value = '  -4.21%  '
$ws.Range('D48').Value = '''0.4559'
$ws.Range('E48').Value = '  -0.53%  '
$ws.Range('D49').Value = '''1.574'
$ws.Range('E49').Value = '  -8.83%  '
$ws.Range('D50').Value = '''0.05499'
$ws.Range('E50').Value = '  -2.85%  '
$ws.Range('D51').Value = '''6.397'
$ws.Range('E51').Value = '  -8.01%  '
